$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row changes ---
# Column Q's header label changes from "Dac_fssel_val" to "tx_gain_row"
$ws.Range("Q1").Value = "tx_gain_row"
# New column AH with header "iftx_stg3_val"
$ws.Range("AH1").Value = "iftx_stg3_val"

# --- Data row updates (rows 2-9) ---
# A column (Temp) changes for rows 6-9: -40 -> -10
$ws.Range("A6").Value = -10
$ws.Range("A7").Value = -10
$ws.Range("A8").Value = -10
$ws.Range("A9").Value = -10

# C column (ChipChannel): all rows 2-9 become 2
$ws.Range("C2:C9").Value = 2

# E column (XIF_1): all rows 2-9 become 1
$ws.Range("E2:E9").Value = 1

# I column (XIF_5): all rows 2-9 become 0
$ws.Range("I2:I9").Value = 0

# L column (XIF_Matrix): rows 2,3,6,7 -> 1 ; rows 5,9 -> 3 ; rows 4,8 unchanged (stay 3)
$ws.Range("L2").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("L5").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 1
$ws.Range("L9").Value = 3

# N column (MSC_table): all rows 2-9 become 1
$ws.Range("N2:N9").Value = 1

# Q column (Dac_fssel_val / now tx_gain_row): rows 2,4,6,8 -> 0 ; rows 3,5,7,9 stay 1
$ws.Range("Q2").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("Q8").Value = 0

# AG column (note) becomes a text cell "Note!!" for all rows 2-9
$ws.Range("AG2:AG9").Value = "Note!!"

# AH column (new, iftx_stg3_val) gets value 1 for rows 2-9
$ws.Range("AH2:AH9").Value = 1

# --- Remove old rows 10-13 (lineup trimmed down to 8 data rows) ---
$ws.Rows("10:13").Delete()
